$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.914.20"
$ws.Range("E2").Value = "  -2.42%  "

$ws.Range("D3").Value = "3.055.01"
$ws.Range("E3").Value = "  -4.62%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.88%  "

$ws.Range("E7").Value = "  -6.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +14.84%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "3.053.66"
$ws.Range("E10").Value = "  -4.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.179"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("E13").Value = "  -7.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("D15").Value = "87.838.39"
$ws.Range("E15").Value = "  -2.02%  "

$ws.Range("D16").Value = "3.627.48"
$ws.Range("E16").Value = "  -4.19%  "

$ws.Range("E17").Value = "  -5.02%  "

$ws.Range("D18").Value = "3.065.31"
$ws.Range("E18").Value = "  -4.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000199"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.94%  "

$ws.Range("E26").Value = "  -2.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").Value = "3.232.00"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.170"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.79%  "

$ws.Range("E32").Value = "  -5.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "503.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.57%  "

$ws.Range("E34").Value = "  -12.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.72%  "

$ws.Range("E36").Value = "  -6.61%  "

$ws.Range("E37").Value = "  -8.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  -4.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("E45").Value = "  -7.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.131"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.59%  "

$ws.Range("E47").Value = "  -2.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0687"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "156.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.18%  "

$ws.Range("E50").Value = "  -5.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.698"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.22%  "
